# ---------------------------------------------------------------------------
# removed trial pages and tests. Added basic pages for sauce demo and a
# sample test; updated config file to point to saucedemo site.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsRunManager = $wb.Worksheets.Item("RUNMANAGER")
$wsData       = $wb.Worksheets.Item("DATA")

# --- RUNMANAGER sheet --------------------------------------------------
# Keep only the header row + the loginLogoutTest row; drop the "newTest"
# and "amazonTest" trial rows. Flip execute flag for the surviving test
# to "yes" now that it points at the real sauce-demo run.
$wsRunManager.Rows.Item(4).Delete()
$wsRunManager.Rows.Item(3).Delete()
$wsRunManager.Range("C2").Value = "yes"

# --- DATA sheet ----------------------------------------------------------
# Drop the trial rows (firefox / newTest / duplicate login rows / amazon),
# keep two sauce-demo data rows, and drop the now-unused fname/menutext
# columns.
$wsData.Rows.Item(7).Delete()
$wsData.Rows.Item(6).Delete()
$wsData.Rows.Item(5).Delete()
$wsData.Rows.Item(4).Delete()

$wsData.Range("B2").Value = "yes"
$wsData.Range("C2").Value = "chrome"
$wsData.Range("D2").Value = "standard_user"
$wsData.Range("E2").Value = "secret_sauce"

$wsData.Range("B3").Value = "yes"
$wsData.Range("C3").Value = "chrome"
$wsData.Range("D3").Value = "locked_out_user"
$wsData.Range("E3").Value = "secret_sauce"

$wsData.Columns.Item(7).Delete()
$wsData.Columns.Item(6).Delete()

# --- Selection / active tab ----------------------------------------------
# RUNMANAGER becomes the active/visible tab with F1:H1048576 selected;
# DATA keeps D3 selected for when it is revisited.
$wsData.Range("D3").Select()
$wsRunManager.Activate()
$wsRunManager.Range("F1:H1048576").Select()
